$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (new most-recent quarter + the one before it).
# Existing data in D:K shifts right to F:M.
$ws.Columns("D:E").Insert()

# Copy number/date formatting from the (now-shifted) F:G columns into the
# freshly inserted D:E columns for the data rows (7 onward) so the new cells
# pick up the same styles (date format / number format) as the rest of the
# table. Rows 5-6 have no data in these columns, so we start at row 7.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @{Row=7; D=43465; E=43373},
    @{Row=8; D=117800; E=118000},
    @{Row=9; D=40600; E=38900},
    @{Row=10; D=77200; E=79100},
    @{Row=12; D=32600; E=29300},
    @{Row=13; D=0; E=0},
    @{Row=14; D=16400; E=0},
    @{Row=15; D=0; E=0},
    @{Row=17; D=109600; E=89200},
    @{Row=18; D=8200; E=28800},
    @{Row=20; D=9900; E=14000},
    @{Row=21; D="NA"; E="NA"},
    @{Row=22; D="NA"; E="NA"},
    @{Row=23; D=18000; E=42900},
    @{Row=24; D=8000; E=-12300},
    @{Row=25; D=0; E=0},
    @{Row=26; D=10100; E=55200},
    @{Row=27; D=10100; E=55200},
    @{Row=28; D=0; E=0},
    @{Row=29; D=0; E=0},
    @{Row=30; D=0; E=0},
    @{Row=31; D=0; E=0},
    @{Row=32; D=-9900; E=-14000},
    @{Row=33; D=10100; E=55200},
    @{Row=34; D=0; E=0},
    @{Row=35; D=10100; E=55200},
    @{Row=38; D=43465; E=43373},
    @{Row=41; D=454500; E=344400},
    @{Row=42; D=190100; E=275700},
    @{Row=43; D=57400; E=56900},
    @{Row=44; D=0; E=0},
    @{Row=45; D=725800; E=685300},
    @{Row=46; D=1427800; E=1362300},
    @{Row=47; D="NA"; E="NA"},
    @{Row=48; D=170700; E=172700},
    @{Row=49; D=24200; E=42300},
    @{Row=50; D=0; E=0},
    @{Row=51; D=0; E=0},
    @{Row=52; D=342800; E=18000},
    @{Row=53; D=0; E=0},
    @{Row=54; D=1965500; E=1595400},
    @{Row=57; D=753100; E=630800},
    @{Row=58; D=0; E=0},
    @{Row=59; D=63600; E=50500},
    @{Row=60; D=816600; E=681200},
    @{Row=61; D=220000; E=0},
    @{Row=62; D=97200; E=97200},
    @{Row=63; D=0; E=0},
    @{Row=64; D=0; E=0},
    @{Row=65; D=0; E=0},
    @{Row=66; D=1135700; E=780400},
    @{Row=68; D=0; E=0},
    @{Row=69; D=0; E=0},
    @{Row=70; D=0; E=0},
    @{Row=71; D=0; E=0},
    @{Row=72; D="NA"; E="NA"},
    @{Row=73; D=0; E=0},
    @{Row=74; D=0; E=0},
    @{Row=75; D=0; E=0},
    @{Row=76; D=829700; E=815000},
    @{Row=77; D=0; E=0},
    @{Row=80; D=43465; E=43373},
    @{Row=81; D=10100; E=55200},
    @{Row=83; D=0; E=0},
    @{Row=84; D=0; E=0},
    @{Row=85; D=0; E=0},
    @{Row=86; D=0; E=0},
    @{Row=87; D=0; E=0},
    @{Row=88; D=0; E=0},
    @{Row=89; D=0; E=0},
    @{Row=91; D=0; E=0},
    @{Row=92; D=0; E=0},
    @{Row=93; D=0; E=0},
    @{Row=94; D=0; E=0},
    @{Row=96; D=0; E=0},
    @{Row=97; D=0; E=0},
    @{Row=98; D=0; E=0},
    @{Row=99; D=0; E=0},
    @{Row=100; D=0; E=0},
    @{Row=101; D=0; E=0},
    @{Row=102; D=0; E=0}
)

foreach ($item in $data) {
    $ws.Range("D" + $item.Row).Value2 = $item.D
    $ws.Range("E" + $item.Row).Value2 = $item.E
}

Write-Output "done"
